$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $v1 = $range1.Value()
    $v2 = $range2.Value()

    $range1.Value = $v2
    $range2.Value = $v1
}

Swap-Rows 26 27
Swap-Rows 160 161
